# Add pH calibration data from 10/2/2019 JFields run (rows 49-51 new data,
# row 52 gets a trailing note). This fills in the previously-blank rows
# 49:52 that only had the Batch value / % off formula present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 49 ----
$ws.Range("A49").Value = 43719
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)   # xlPasteFormats - copy date number format
$ws.Range("B49").Value = 2217.19
$ws.Range("D49").Formula = "=100*(B49-C49)/C49"
$ws.Range("E49").Value = 169
$ws.Range("F49").Value = "run at beginning of day with new probe"

# ---- Row 50 ----
$ws.Range("A50").Value = 43720
$ws.Range("A48").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("B50").Value = 2218.23
$ws.Range("D50").Formula = "=100*(B50-C50)/C50"
$ws.Range("E50").Value = 169
$ws.Range("F50").Value = "run at beginning of day"

# ---- Row 51 ----
$ws.Range("A51").Value = 43723
$ws.Range("A48").Copy()
$ws.Range("A51").PasteSpecial(-4122)
$ws.Range("B51").Value = 2091.32236590917
$ws.Range("D51").Formula = "=100*(B51-C51)/C51"
$ws.Range("E51").Value = 169
$ws.Range("F51").Value = "opened crm (9/8/2019); ph calibration off?"
$ws.Range("A48").Copy()
$ws.Range("F51").PasteSpecial(-4122)

# ---- Row 52 ----
$ws.Range("F52").Value = "opened crm (8/7/2019)"

# Recalculate the shared "% off" formula column so cached <v> values match
$excel.Calculate()

# Update the window selection to the last-edited cell
$ws.Range("F52").Select()
